# Auto-generated edit script: refresh rolling 169-hour prediction window
# (shift the Date/Interval/Lookup sequence forward by 5 hours and write new
# retrained model Prediction values), per "Retraining models with the latest data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates   = @(45975,45975,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45976,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45977,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45978,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45979,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45980,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45981,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982,45982)
$hours   = @(23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)
$preds   = @(0.104,0.104,0.01,0.01,0.0,0.015,0.017,0.017,0.014,0.0,0.022,0.27,0.725,1.144,1.531,1.523,1.447,1.02,0.591,0.113,0.0,0.0,0.01,0.0,0.026,0.028,0.027,0.027,0.03,0.027,0.027,0.029,0.0,0.0,0.043,0.407,0.948,1.646,1.99,2.181,1.98,1.525,0.646,0.121,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.0,0.01,0.015,0.0,0.01,0.017,0.017,0.041,0.31,0.86,1.544,1.542,1.587,1.165,0.752,0.296,0.043,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.023,0.136,0.414,0.639,0.775,0.777,0.635,0.374,0.15,0.033,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.019,0.133,0.371,0.627,0.721,0.737,0.64,0.451,0.212,0.033,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.014,0.017,0.017,0.017,0.017,0.02,0.216,0.62,0.925,1.351,1.389,1.178,0.783,0.371,0.048,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.017,0.131,0.422,0.703,0.775,0.791,0.695,0.456,0.203,0.034,0.0,0.0,0.0,0.0,0.0)
$lookups = @("14.11.202523","14.11.202524","15.11.20251","15.11.20252","15.11.20253","15.11.20254","15.11.20255","15.11.20256","15.11.20257","15.11.20258","15.11.20259","15.11.202510","15.11.202511","15.11.202512","15.11.202513","15.11.202514","15.11.202515","15.11.202516","15.11.202517","15.11.202518","15.11.202519","15.11.202520","15.11.202521","15.11.202522","15.11.202523","15.11.202524","16.11.20251","16.11.20252","16.11.20253","16.11.20254","16.11.20255","16.11.20256","16.11.20257","16.11.20258","16.11.20259","16.11.202510","16.11.202511","16.11.202512","16.11.202513","16.11.202514","16.11.202515","16.11.202516","16.11.202517","16.11.202518","16.11.202519","16.11.202520","16.11.202521","16.11.202522","16.11.202523","16.11.202524","17.11.20251","17.11.20252","17.11.20253","17.11.20254","17.11.20255","17.11.20256","17.11.20257","17.11.20258","17.11.20259","17.11.202510","17.11.202511","17.11.202512","17.11.202513","17.11.202514","17.11.202515","17.11.202516","17.11.202517","17.11.202518","17.11.202519","17.11.202520","17.11.202521","17.11.202522","17.11.202523","17.11.202524","18.11.20251","18.11.20252","18.11.20253","18.11.20254","18.11.20255","18.11.20256","18.11.20257","18.11.20258","18.11.20259","18.11.202510","18.11.202511","18.11.202512","18.11.202513","18.11.202514","18.11.202515","18.11.202516","18.11.202517","18.11.202518","18.11.202519","18.11.202520","18.11.202521","18.11.202522","18.11.202523","18.11.202524","19.11.20251","19.11.20252","19.11.20253","19.11.20254","19.11.20255","19.11.20256","19.11.20257","19.11.20258","19.11.20259","19.11.202510","19.11.202511","19.11.202512","19.11.202513","19.11.202514","19.11.202515","19.11.202516","19.11.202517","19.11.202518","19.11.202519","19.11.202520","19.11.202521","19.11.202522","19.11.202523","19.11.202524","20.11.20251","20.11.20252","20.11.20253","20.11.20254","20.11.20255","20.11.20256","20.11.20257","20.11.20258","20.11.20259","20.11.202510","20.11.202511","20.11.202512","20.11.202513","20.11.202514","20.11.202515","20.11.202516","20.11.202517","20.11.202518","20.11.202519","20.11.202520","20.11.202521","20.11.202522","20.11.202523","20.11.202524","21.11.20251","21.11.20252","21.11.20253","21.11.20254","21.11.20255","21.11.20256","21.11.20257","21.11.20258","21.11.20259","21.11.202510","21.11.202511","21.11.202512","21.11.202513","21.11.202514","21.11.202515","21.11.202516","21.11.202517","21.11.202518","21.11.202519","21.11.202520","21.11.202521","21.11.202522","21.11.202523")

$n = $dates.Length
for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $hours[$i]
    $ws.Cells.Item($r, 3).Value = $preds[$i]
    $ws.Cells.Item($r, 4).Value = $lookups[$i]
}
